# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to reflect newly generated numbers from the site build at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F3").Value = 532
$sheet1.Range("F8").Value = 2283
$sheet1.Range("F10").Value = 5677
$sheet1.Range("F12").Value = 372

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F4").Value = 532
$sheet4.Range("F11").Value = 2283
$sheet4.Range("F13").Value = 5677
$sheet4.Range("F15").Value = 372
